$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the "float" description cells.
#   H2 gains a unit annotation; J2 reverts to the plain "#float" label.
$ws.Range("H2").Value = "#float,  unit:mlormg"
$ws.Range("J2").Value = "#float"

# Row 3: new enum/description row mapping each header to a French keyword.
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
